$wb = $excel.ActiveWorkbook
$wsInstructions = $wb.Worksheets.Item("Instructions")
$ws = $wb.Worksheets.Item("InflowWind")

# ---------------------------------------------------------------------------
# Add the new "Wind Sensor Measurements" category (lidar) to the InflowWind
# OutList worksheet: one category-header row (30) followed by five data rows
# (31-35) describing WindMeas1..WindMeas5.
# ---------------------------------------------------------------------------

# Clear out the old blank placeholder rows 30-32 (they only carried leftover
# formatting in columns E/F) so we start from a clean slate for rows 30-35.
$ws.Range("A30:G35").ClearFormats() | Out-Null
$ws.Range("A30:G35").ClearContents() | Out-Null

# Row 30: category header, formatted like row 2 ("Wind Motions").
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:E2").Copy() | Out-Null
$ws.Range("C30:E30").PasteSpecial(-4122) | Out-Null
$ws.Range("A30").Value = "Wind Sensor Measurements"

# Rows 31-35: data rows, formatted like row 3 (a normal OutList data row).
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B31:B35").PasteSpecial(-4122) | Out-Null
$ws.Range("D3:F3").Copy() | Out-Null
$ws.Range("D31:F31").PasteSpecial(-4122) | Out-Null
$ws.Range("D32:F32").PasteSpecial(-4122) | Out-Null
$ws.Range("D33:F33").PasteSpecial(-4122) | Out-Null
$ws.Range("D34:F34").PasteSpecial(-4122) | Out-Null
$ws.Range("D35:F35").PasteSpecial(-4122) | Out-Null

$ws.Range("B31").Value = "WindMeas1"
$ws.Range("B32").Value = "WindMeas2"
$ws.Range("B33").Value = "WindMeas3"
$ws.Range("B34").Value = "WindMeas4"
$ws.Range("B35").Value = "WindMeas5"

$ws.Range("D31").Value = "Wind measurement at sensor 1"
$ws.Range("D32").Value = "Wind measurement at sensor 2"
$ws.Range("D33").Value = "Wind measurement at sensor 3"
$ws.Range("D34").Value = "Wind measurement at sensor 4"
$ws.Range("D35").Value = "Wind measurement at sensor 5"

$ws.Range("E31").Value = "Defined by sensor"
$ws.Range("E32").Value = "Defined by sensor"
$ws.Range("E33").Value = "Defined by sensor"
$ws.Range("E34").Value = "Defined by sensor"
$ws.Range("E35").Value = "Defined by sensor"

$ws.Range("F31").Value = "(m/s)"
$ws.Range("F32").Value = "(m/s)"
$ws.Range("F33").Value = "(m/s)"
$ws.Range("F34").Value = "(m/s)"
$ws.Range("F35").Value = "(m/s)"

$ws.Range("G31").Value = "p%lidar%SensorType == SensorType_None"
$ws.Range("G32").Value = "p%lidar%NumPulseGate < 2"
$ws.Range("G33").Value = "p%lidar%NumPulseGate < 3"
$ws.Range("G34").Value = "p%lidar%NumPulseGate < 4"
$ws.Range("G35").Value = "p%lidar%NumPulseGate < 5"

# ---------------------------------------------------------------------------
# The Instructions sheet's D8 cell (=COUNTA(InflowWind!B2:B1333)) recalculates
# automatically from 27 to 32 now that five more names exist in column B.
# ---------------------------------------------------------------------------

# Restore the selections recorded in the saved workbook: Instructions!D8 and
# InflowWind!B31 (InflowWind remains the active/tab-selected sheet).
$wsInstructions.Range("D8").Select() | Out-Null
$ws.Range("B31").Select() | Out-Null
